# update of the branch
#
# 1) The "update automatically" date footer field (type="datetimeFigureOut")
#    on the slide master and on every slide layout had its cached display
#    text bumped from "06-Apr-22" to "08-Aug-22".
# 2) The picture placeholder on slide 4 ("Filter on column names") was
#    moved/resized.

$p = $ppt.ActivePresentation

$oldDate = "06-Apr-22"
$newDate = "08-Aug-22"

# ppPlaceholderDate
$ppPlaceholderDate = 16

$msoPlaceholder = 14

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Type -eq $msoPlaceholder -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every layout off the master has its own (inherited/overridden) date
# placeholder shape that also caches the rendered text.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Reposition / resize the picture ("Content Placeholder 4") on slide 4.
$slide4 = $p.Slides.Item(4)
$picShape = $slide4.Shapes.Item("Content Placeholder 4")
$picShape.Left = 205.8473
$picShape.Top = 196.8886
$picShape.Width = 525.82663
$picShape.Height = 163.89403
